$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update E14 value (figures fix)
$ws.Range("E14").Value = 112.842

# Move the active selection to E15 (as recorded in the saved sheet view)
$ws.Range("E15").Select()
